$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-looking-like-numbers are stored as text (preserve exact formatting,
# e.g. trailing zeros / avoid float rounding) by forcing Text number format
# before assigning their values.
$textCells = @("D5","D8","D16","D18","D22","D25","D34","D37","D39","D40","D41","D43","D45","D46","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "27.920.20"

# Row 3
$ws.Range("D3").Value = "1.643.17"
$ws.Range("E3").Value = "  +1.22%  "

# Row 5
$ws.Range("D5").Value = "213.58"
$ws.Range("E5").Value = "  +0.91%  "

# Row 6
$ws.Range("E6").Value = "  -0.14%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").Value = "23.62"
$ws.Range("E8").Value = "  +1.60%  "

# Row 9
$ws.Range("E9").Value = "  +0.53%  "

# Row 11
$ws.Range("E11").Value = "  -1.92%  "

# Row 12
$ws.Range("D12").Value = "1.875.41"
$ws.Range("E12").Value = "  +1.19%  "

# Row 13
$ws.Range("D13").Value = "1.642.46"
$ws.Range("E13").Value = "  +1.85%  "

# Row 14
$ws.Range("E14").Value = "  +4.42%  "

# Row 15
$ws.Range("E15").Value = "  +0.65%  "

# Row 16
$ws.Range("D16").Value = "65.79"
$ws.Range("E16").Value = "  +0.63%  "

# Row 17
$ws.Range("D17").Value = "27.898.19"
$ws.Range("E17").Value = "  +1.34%  "

# Row 18
$ws.Range("D18").Value = "230.41"
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("E19").Value = "  +0.86%  "

# Row 20
$ws.Range("E20").Value = "  +1.09%  "

# Row 21
$ws.Range("E21").Value = "  +0.00%  "

# Row 22
$ws.Range("D22").Value = "10.91"
$ws.Range("E22").Value = "  +4.58%  "

# Row 23
$ws.Range("E23").Value = "  +1.46%  "

# Row 24
$ws.Range("E24").Value = "  +2.47%  "

# Row 25
$ws.Range("D25").Value = "152.13"
$ws.Range("E25").Value = "  +1.72%  "

# Row 26
$ws.Range("E26").Value = "  +0.60%  "

# Row 27
$ws.Range("E27").Value = "  +0.75%  "

# Row 28
$ws.Range("E28").Value = "  +1.08%  "

# Row 29
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("E30").Value = "  +0.96%  "

# Row 32
$ws.Range("E32").Value = "  +2.02%  "

# Row 33
$ws.Range("D33").Value = "1.425.53"
$ws.Range("E33").Value = "  -2.54%  "

# Row 34
$ws.Range("D34").Value = "3.10"
$ws.Range("E34").Value = "  +0.85%  "

# Row 35
$ws.Range("E35").Value = "  +1.64%  "

# Row 36
$ws.Range("E36").Value = "  -0.25%  "

# Row 37
$ws.Range("D37").Value = "0.890"
$ws.Range("E37").Value = "  +1.83%  "

# Row 38
$ws.Range("E38").Value = "  +0.61%  "

# Row 39 (was TrustWalletToken, becomes ImmutableX)
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.557"
$ws.Range("E39").Value = "  +0.42%  "

# Row 40 (was ImmutableX, becomes TrustWalletToken)
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.920"
$ws.Range("E40").Value = "  -2.53%  "

# Row 41
$ws.Range("D41").Value = "1.03"
$ws.Range("E41").Value = "  +2.59%  "

# Row 43
$ws.Range("D43").Value = "68.56"
$ws.Range("E43").Value = "  +1.41%  "

# Row 44
$ws.Range("E44").Value = "  +1.10%  "

# Row 45
$ws.Range("D45").Value = "5.45"
$ws.Range("E45").Value = "  +2.87%  "

# Row 46
$ws.Range("D46").Value = "1.81"
$ws.Range("E46").Value = "  +2.64%  "

# Row 47
$ws.Range("E47").Value = "  +0.21%  "

# Row 48
$ws.Range("D48").Value = "1.784.16"

# Row 49
$ws.Range("D49").Value = "89.10"
$ws.Range("E49").Value = "  +1.98%  "

# Row 50
$ws.Range("E50").Value = "  +0.47%  "

# Row 51
$ws.Range("E51").Value = "  +0.64%  "
